$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.449328780174255
$ws.Range("B1").Value = 3.611166477203369
$ws.Range("C1").Value = 6.244411945343018
$ws.Range("D1").Value = 1.519080877304077
$ws.Range("E1").Value = 0.8900896906852722
